# sc5 pro power on success
# - Clear the Mode/AF/Pull/Output-type/Output-speed (D:H) columns on a set of
#   power/ground/special pins on the "pin" sheet (they don't apply to these pins).
# - Rework the last 4 rows of the "power" sheet's reset / DDR power-good block
#   into a 3-row Name/Type table driven in FUNCTION mode (drop the 4th row).
# - Restore the saved scroll position / selection on both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "pin" sheet
# ---------------------------------------------------------------------------
$pin = $wb.Worksheets.Item("pin")

$rowsToClear = @(7, 11, 15, 20, 21, 22, 23, 28, 29, 50, 51, 73, 74, 75, 76, 77, 88, 95, 100, 101)
foreach ($r in $rowsToClear) {
    $pin.Range("D" + $r + ":H" + $r).ClearContents()
}

# ---------------------------------------------------------------------------
# "power" sheet
# ---------------------------------------------------------------------------
$power = $wb.Worksheets.Item("power")

# Row 16: SYS-RST-DEASSERT / SYS_RST_DEASSERT, FUNCTION mode, 0 / 1000
$power.Range("A16").Value = "SYS-RST-DEASSERT"
$power.Range("B16").Value = "SYS_RST_DEASSERT"
$power.Range("C16:E16").ClearContents()
$power.Range("F16").Value = "FUNCTION"
$power.Range("G16").Value = 0
$power.Range("H16").Value = 1000

# Row 17: PG-DDR-0 / DDR_POWER_GOOD, FUNCTION mode, 0 / 30000
$power.Range("A17").Value = "PG-DDR-0"
$power.Range("B17").Value = "DDR_POWER_GOOD"
$power.Range("C17:E17").ClearContents()
$power.Range("F17").Value = "FUNCTION"
$power.Range("G17").Value = 0
$power.Range("H17").Value = 30000

# Row 18: SYS-RST-ASSERT / SYS_RST_ASSERT, FUNCTION mode, 0 / 30000
$power.Range("A18").Value = "SYS-RST-ASSERT"
$power.Range("B18").Value = "SYS_RST_ASSERT"
$power.Range("C18:E18").ClearContents()
$power.Range("F18").Value = "FUNCTION"
$power.Range("G18").Value = 0
$power.Range("H18").Value = 30000

# Row 19 (old PG-DDR-1 row) is no longer needed - delete it entirely.
$power.Rows.Item(19).Delete()

# ---------------------------------------------------------------------------
# View state: scroll position + selection on both sheets
# ---------------------------------------------------------------------------
$pin.Activate()
$pinWin = $excel.ActiveWindow
$pin.Range("A31").Select()
$pinWin.ScrollRow = 31
$pinWin.ScrollColumn = 1
$pin.Range("F48").Select()

$power.Activate()
$powerWin = $excel.ActiveWindow
$power.Range("A1").Select()
$powerWin.ScrollRow = 1
$powerWin.ScrollColumn = 1
$power.Range("C18").Select()

$pin.Activate()
